# feat: adds default buckets
# Replace the placeholder "some_bucket2" s3-bucket values in rows 3 and 4
# with the new default bucket names "open" and "scratch".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "open"
$ws.Range("E4").Value = "scratch"

# Leave the selection where the author's last save left it.
[void]$ws.Range("E5").Select()
